# Update the NATMI Angpt4-Tek LR-pair sheet with the freshly recomputed
# TPM-based statistics. Rows 2-9 get new numeric values (and a few of the
# "Target cluster" / "Sending cluster" labels are reshuffled since the
# underlying row ordering changed), and the old rows 10-11 (which paired
# MuSCs -> MuSCs / MuSCs -> Resolving-Mac) are dropped entirely because
# the "Resolving-Mac" cluster no longer appears in the refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Angpt4"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5282566666666667
$ws.Range("H2").Value = 1.58477
$ws.Range("I2").Value = 0.8498693102778859
$ws.Range("J2").Value = 0.8498693102778859
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 71.47459166666665
$ws.Range("N2").Value = 214.423775
$ws.Range("O2").Value = 0.9365945886638486
$ws.Range("P2").Value = 0.9365945886638485
$ws.Range("Q2").Value = 37.75692954519444
$ws.Range("R2").Value = 339.81236590675
$ws.Range("S2").Value = 0.7959829970777452
$ws.Range("T2").Value = 0.7959829970777452
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Angpt4"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5282566666666667
$ws.Range("H3").Value = 1.58477
$ws.Range("I3").Value = 0.8498693102778859
$ws.Range("J3").Value = 0.8498693102778859
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.258629666666666
$ws.Range("N3").Value = 9.775889
$ws.Range("O3").Value = 0.04270069742396077
$ws.Range("P3").Value = 0.04270069742396076
$ws.Range("Q3").Value = 1.721392845614445
$ws.Range("R3").Value = 15.49253561053
$ws.Range("S3").Value = 0.03629001226808624
$ws.Range("T3").Value = 0.03629001226808623
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Angpt4"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5282566666666667
$ws.Range("H4").Value = 1.58477
$ws.Range("I4").Value = 0.8498693102778859
$ws.Range("J4").Value = 0.8498693102778859
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1062743333333333
$ws.Range("N4").Value = 0.318823
$ws.Range("O4").Value = 0.001392606284175224
$ws.Range("P4").Value = 0.001392606284175224
$ws.Range("Q4").Value = 0.0561401250788889
$ws.Range("R4").Value = 0.50526112571
$ws.Range("S4").Value = 0.001183533342220647
$ws.Range("T4").Value = 0.001183533342220647
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Angpt4"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5282566666666667
$ws.Range("H5").Value = 1.58477
$ws.Range("I5").Value = 0.8498693102778859
$ws.Range("J5").Value = 0.8498693102778859
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.47377
$ws.Range("N5").Value = 4.42131
$ws.Range("O5").Value = 0.01931210762801542
$ws.Range("P5").Value = 0.01931210762801541
$ws.Range("Q5").Value = 0.7785288276333334
$ws.Range("R5").Value = 7.0067594487
$ws.Range("S5").Value = 0.01641276758983376
$ws.Range("T5").Value = 0.01641276758983376
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Angpt4"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.09331733333333335
$ws.Range("H6").Value = 0.279952
$ws.Range("I6").Value = 0.1501306897221141
$ws.Range("J6").Value = 0.1501306897221141
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 71.47459166666665
$ws.Range("N6").Value = 214.423775
$ws.Range("O6").Value = 0.9365945886638486
$ws.Range("P6").Value = 0.9365945886638485
$ws.Range("Q6").Value = 6.669818295422222
$ws.Range("R6").Value = 60.0283646588
$ws.Range("S6").Value = 0.1406115915861033
$ws.Range("T6").Value = 0.1406115915861033
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Angpt4"
$ws.Range("C7").Value = "Tek"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.09331733333333335
$ws.Range("H7").Value = 0.279952
$ws.Range("I7").Value = 0.1501306897221141
$ws.Range("J7").Value = 0.1501306897221141
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.258629666666666
$ws.Range("N7").Value = 9.775889
$ws.Range("O7").Value = 0.04270069742396077
$ws.Range("P7").Value = 0.04270069742396076
$ws.Range("Q7").Value = 0.3040866308142222
$ws.Range("R7").Value = 2.736779677328
$ws.Range("S7").Value = 0.00641068515587453
$ws.Range("T7").Value = 0.00641068515587453
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Angpt4"
$ws.Range("C8").Value = "Tek"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.09331733333333335
$ws.Range("H8").Value = 0.279952
$ws.Range("I8").Value = 0.1501306897221141
$ws.Range("J8").Value = 0.1501306897221141
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1062743333333333
$ws.Range("N8").Value = 0.318823
$ws.Range("O8").Value = 0.001392606284175224
$ws.Range("P8").Value = 0.001392606284175224
$ws.Range("Q8").Value = 0.009917237388444447
$ws.Range("R8").Value = 0.08925513649600002
$ws.Range("S8").Value = 0.0002090729419545768
$ws.Range("T8").Value = 0.0002090729419545768
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Angpt4"
$ws.Range("C9").Value = "Tek"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.09331733333333335
$ws.Range("H9").Value = 0.279952
$ws.Range("I9").Value = 0.1501306897221141
$ws.Range("J9").Value = 0.1501306897221141
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.47377
$ws.Range("N9").Value = 4.42131
$ws.Range("O9").Value = 0.01931210762801542
$ws.Range("P9").Value = 0.01931210762801541
$ws.Range("Q9").Value = 0.1375282863466667
$ws.Range("R9").Value = 1.23775457712
$ws.Range("S9").Value = 0.002899340038181655
$ws.Range("T9").Value = 0.002899340038181655

# Remove rows 10 and 11 (no longer present in the data)
$ws.Rows("10:11").Delete()
